# Auto update Excel log - append new sensor readings to PIR, Humidity, and Temperature sheets
$wb = $excel.ActiveWorkbook

# ---- PIR sheet ----
$wsPIR = $wb.Worksheets.Item("PIR")
$PIRData = @(
    @("2026-01-28", "14:58:49", "14:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-28", "14:58:50", "14:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-28", "14:58:55", "14:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-28", "14:58:59", "14:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-28", "14:59:04", "14:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-28", "14:59:09", "14:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-28", "14:59:15", "14:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-28", "14:59:19", "14:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-28", "14:59:24", "14:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-28", "14:59:29", "14:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-28", "14:59:35", "14:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-28", "14:59:39", "14:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-28", "14:59:44", "14:00", "Bathroom", "No Motion", "Inactive")
)
$startRow = 135
# Keep all columns as plain text so Excel does not auto-convert
# strings like 2026-01-28 or 87.2% into dates/numbers.
$wsPIR.Range("A135:F147").NumberFormat = "@"
for ($i = 0; $i -lt $PIRData.Count; $i++) {
    $r = $startRow + $i
    $row = $PIRData[$i]
    $wsPIR.Cells.Item($r, 1).Value = $row[0]
    $wsPIR.Cells.Item($r, 2).Value = $row[1]
    $wsPIR.Cells.Item($r, 3).Value = $row[2]
    $wsPIR.Cells.Item($r, 4).Value = $row[3]
    $wsPIR.Cells.Item($r, 5).Value = $row[4]
    $wsPIR.Cells.Item($r, 6).Value = $row[5]
}

# ---- Humidity sheet ----
$wsHum = $wb.Worksheets.Item("Humidity")
$HumData = @(
    @("2026-01-28", "14:58:49", "14:00", "Bathroom", "87.2%", "Active"),
    @("2026-01-28", "14:58:51", "14:00", "Bathroom", "88.1%", "Active"),
    @("2026-01-28", "14:58:54", "14:00", "Bathroom", "87.2%", "Active"),
    @("2026-01-28", "14:59:02", "14:00", "Bathroom", "88.1%", "Active"),
    @("2026-01-28", "14:59:06", "14:00", "Bathroom", "88.1%", "Active"),
    @("2026-01-28", "14:59:10", "14:00", "Bathroom", "88.1%", "Active"),
    @("2026-01-28", "14:59:14", "14:00", "Bathroom", "87.2%", "Active"),
    @("2026-01-28", "14:59:22", "14:00", "Bathroom", "88.2%", "Active"),
    @("2026-01-28", "14:59:26", "14:00", "Bathroom", "87.2%", "Active"),
    @("2026-01-28", "14:59:34", "14:00", "Bathroom", "87.2%", "Active"),
    @("2026-01-28", "14:59:38", "14:00", "Bathroom", "88.2%", "Active"),
    @("2026-01-28", "14:59:42", "14:00", "Bathroom", "88.2%", "Active"),
    @("2026-01-28", "14:59:46", "14:00", "Bathroom", "87.2%", "Active")
)
$startRow = 131
# Keep all columns as plain text so Excel does not auto-convert
# strings like 2026-01-28 or 87.2% into dates/numbers.
$wsHum.Range("A131:F143").NumberFormat = "@"
for ($i = 0; $i -lt $HumData.Count; $i++) {
    $r = $startRow + $i
    $row = $HumData[$i]
    $wsHum.Cells.Item($r, 1).Value = $row[0]
    $wsHum.Cells.Item($r, 2).Value = $row[1]
    $wsHum.Cells.Item($r, 3).Value = $row[2]
    $wsHum.Cells.Item($r, 4).Value = $row[3]
    $wsHum.Cells.Item($r, 5).Value = $row[4]
    $wsHum.Cells.Item($r, 6).Value = $row[5]
}

# ---- Temperature sheet ----
$wsTemp = $wb.Worksheets.Item("Temperature")
$TempData = @(
    @("2026-01-28", "14:58:50", "14:00", "Bathroom", "22.9C", "Active"),
    @("2026-01-28", "14:58:51", "14:00", "Bathroom", "22.9C", "Active"),
    @("2026-01-28", "14:58:54", "14:00", "Bathroom", "22.9C", "Active"),
    @("2026-01-28", "14:59:02", "14:00", "Bathroom", "22.9C", "Active"),
    @("2026-01-28", "14:59:06", "14:00", "Bathroom", "22.9C", "Active"),
    @("2026-01-28", "14:59:10", "14:00", "Bathroom", "22.9C", "Active"),
    @("2026-01-28", "14:59:14", "14:00", "Bathroom", "22.9C", "Active"),
    @("2026-01-28", "14:59:22", "14:00", "Bathroom", "22.9C", "Active"),
    @("2026-01-28", "14:59:26", "14:00", "Bathroom", "22.8C", "Active"),
    @("2026-01-28", "14:59:34", "14:00", "Bathroom", "22.8C", "Active"),
    @("2026-01-28", "14:59:38", "14:00", "Bathroom", "22.9C", "Active"),
    @("2026-01-28", "14:59:42", "14:00", "Bathroom", "22.9C", "Active"),
    @("2026-01-28", "14:59:46", "14:00", "Bathroom", "22.9C", "Active")
)
$startRow = 131
# Keep all columns as plain text so Excel does not auto-convert
# strings like 2026-01-28 or 87.2% into dates/numbers.
$wsTemp.Range("A131:F143").NumberFormat = "@"
for ($i = 0; $i -lt $TempData.Count; $i++) {
    $r = $startRow + $i
    $row = $TempData[$i]
    $wsTemp.Cells.Item($r, 1).Value = $row[0]
    $wsTemp.Cells.Item($r, 2).Value = $row[1]
    $wsTemp.Cells.Item($r, 3).Value = $row[2]
    $wsTemp.Cells.Item($r, 4).Value = $row[3]
    $wsTemp.Cells.Item($r, 5).Value = $row[4]
    $wsTemp.Cells.Item($r, 6).Value = $row[5]
}

